# Respostas.xlsx edit
#  1. Insert a new "Esg" worksheet between "Técnica" and "Consolidado",
#     populated with its header row + 3 response rows.
#  2. Append a new response row (row 11) to the "Técnica" sheet.
#  3. On the "Consolidado" sheet: tweak the K1 header text and zero out a
#     block of previously-non-zero values in rows 2-4.
# The "Comercial" sheet is left untouched (it just shifts from sheet3.xml
# to sheet4.xml on disk, a natural side effect of inserting the new sheet).
#
# NOTE: worksheet object references returned by this host are
# position-bound, not identity-bound. Once a sheet is inserted/moved, any
# previously-fetched handle for a sheet that used to sit at that index can
# silently start pointing at the *new* sheet instead. So every worksheet
# handle below is (re)fetched by name immediately before it is used.

$wb = $excel.ActiveWorkbook

# --- 1. Add the "Esg" sheet right before "Consolidado" -----------------
$esg = $wb.Worksheets.Add($wb.Worksheets.Item("Consolidado"))
$esg.Name = "Esg"

$esg = $wb.Worksheets.Item("Esg")
$esg.Range("A1").Value = "Data"
$esg.Range("B1").Value = "Hora"
$esg.Range("C1").Value = "E-mail"
$esg.Range("D1").Value = "Categoria"
$esg.Range("E1").Value = "Fornecedor"
$esg.Range("F1").Value = "Clima:`nComo o fornecedor está gerenciando as emissões de carbono (medindo, reduzindo e mitigando)"
$esg.Range("G1").Value = "Social:`nComo o fornecedor está comprometido em promover direitos humanos, diversidade e desenvolvimento social em sua cadeia de valor?"
$esg.Range("H1").Value = "Governança:`nO fornecedor tem um código de conduta ou política sobre ética e transparência"
$esg.Range("I1").Value = "Clima:`nComo o fornecedor está gerenciando as emissões de carbono (medindo, reduzindo e mitigando) (PONDERADA)"
$esg.Range("J1").Value = "Social:`nComo o fornecedor está comprometido em promover direitos humanos, diversidade e desenvolvimento social em sua cadeia de valor? (PONDERADA)"
$esg.Range("K1").Value = "Governança:`nO fornecedor tem um código de conduta ou política sobre ética e transparência (PONDERADA)"

$esg = $wb.Worksheets.Item("Esg")
$esg.Range("A2").NumberFormat = "@"
$esg.Range("A2").Value = "23/06/2025"
$esg.Range("A2").ClearFormats()
$esg.Range("B2").Value = "12:23:24"
$esg.Range("C2").Value = "Teste3"
$esg.Range("D2").Value = "3PL"
$esg.Range("E2").Value = "UNIDOCK’S ASSESSORIA E LOGÍSTICA DE MATERIAIS LTDA"
$esg.Range("F2").Value = 1
$esg.Range("G2").Value = 1
$esg.Range("H2").Value = 1
$esg.Range("I2").Value = 0.05
$esg.Range("J2").Value = 0.05
$esg.Range("K2").Value = 0.05

$esg = $wb.Worksheets.Item("Esg")
$esg.Range("A3").NumberFormat = "@"
$esg.Range("A3").Value = "23/06/2025"
$esg.Range("A3").ClearFormats()
$esg.Range("B3").Value = "12:23:33"
$esg.Range("C3").Value = "Teste3"
$esg.Range("D3").Value = "3PL"
$esg.Range("E3").Value = "CEVA LOGISTICS LTDA"
$esg.Range("F3").Value = 2
$esg.Range("G3").Value = 2
$esg.Range("H3").Value = 2
$esg.Range("I3").Value = 0.1
$esg.Range("J3").Value = 0.1
$esg.Range("K3").Value = 0.1

$esg = $wb.Worksheets.Item("Esg")
$esg.Range("A4").NumberFormat = "@"
$esg.Range("A4").Value = "23/06/2025"
$esg.Range("A4").ClearFormats()
$esg.Range("B4").Value = "12:23:40"
$esg.Range("C4").Value = "Teste3"
$esg.Range("D4").Value = "3PL"
$esg.Range("E4").Value = "KUEHNE + NAGEL SERVIÇOS LOGISTICOS LTDA"
$esg.Range("F4").Value = 3
$esg.Range("G4").Value = 3
$esg.Range("H4").Value = 3
$esg.Range("I4").Value = 0.15
$esg.Range("J4").Value = 0.15
$esg.Range("K4").Value = 0.15

# --- 2. Append row 11 to the "Técnica" sheet ----------------------------
$tecnica = $wb.Worksheets.Item("Técnica")

$tecnica.Cells.Item(11, 1).NumberFormat = "@"
$tecnica.Cells.Item(11, 1).Value = "07/07/2025"
$tecnica.Cells.Item(11, 1).ClearFormats()
$tecnica.Cells.Item(11, 2).Value = "11:06:18"
$tecnica.Cells.Item(11, 3).Value = "Teste2"
$tecnica.Cells.Item(11, 4).Value = "CALL CENTER"
$tecnica.Cells.Item(11, 5).Value = "CTX MLB"
$tecnica.Cells.Item(11, 6).Value = 2
$tecnica.Cells.Item(11, 7).Value = 2
$tecnica.Cells.Item(11, 8).Value = 2
$tecnica.Cells.Item(11, 9).Value = 2
$tecnica.Cells.Item(11, 10).Value = 2
$tecnica.Cells.Item(11, 11).Value = 2
$tecnica.Cells.Item(11, 12).Value = 0.4
$tecnica.Cells.Item(11, 13).Value = 0.3
$tecnica.Cells.Item(11, 14).Value = 0.3
$tecnica.Cells.Item(11, 15).Value = 0.2
$tecnica.Cells.Item(11, 16).Value = 0.3
$tecnica.Cells.Item(11, 17).Value = 0.2

# --- 3. Update the "Consolidado" sheet ----------------------------------
$consolidado = $wb.Worksheets.Item("Consolidado")
$consolidado.Range("K1").Value = "Iniciativas de redução de custos:`nO fornecedor demonstrou esforços para reduzir o preço ao propor iniciativas de redução de custos.1"

foreach ($row in 2..4) {
    $consolidado = $wb.Worksheets.Item("Consolidado")
    foreach ($col in @("E", "I", "J", "K", "L", "M", "N", "R")) {
        $consolidado.Range("$col$row").Value = 0
    }
}
